# 2024 08 13 추가 수정
# Fixes / updates to the db설계 (DB design) workbook:
#  - TBL_CLASS.PROGRAMNO column description was mislabeled "강의실번호"
#    (room number) -> corrected to "프로그램번호" (program number)
#  - TBL_PROGRAM.PROGRAM_NAME column type shrunk from VARCHAR(50) to VARCHAR(20)
#  - TBL_PROGRAM gained a new column PROGRAM_ENAME (VARCHAR(20), "프로그램 영문명")
#    replacing the old unused DESCRIPTION/"프로그램 설명" row
#  - Selection cursor left on M18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TBL_CLASS table (rows 5-13): PROGRAMNO row description fix
$ws.Range("H6").Value = "프로그램번호"

# TBL_PROGRAM table (rows 18-20)
# Row 19: PROGRAM_NAME column - type VARCHAR(50) -> VARCHAR(20)
$ws.Range("N19").Value = "VARCHAR(20)"

# Row 20: DESCRIPTION/TEXT/"프로그램 설명" -> PROGRAM_ENAME/VARCHAR(20)/"프로그램 영문명"
$ws.Range("M20").Value = "PROGRAM_ENAME"
$ws.Range("N20").Value = "VARCHAR(20)"
$ws.Range("P20").Value = "프로그램 영문명"

# Column P widened to fit the new longer description text
$ws.Columns.Item(16).ColumnWidth = 15.15

# Leave the selection on the newly edited row
$ws.Range("M18").Select()
